$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.126.06'
$ws.Range('E2').Value = '  -1.40%  '
$ws.Range('D3').Value = '1.657.49'
$ws.Range('E3').Value = '  -1.17%  '
$ws.Range('E4').Value = '  +0.23%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '216.70'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -1.38%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.5165'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -2.75%  '
$ws.Range('E7').Value = '  +0.22%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.2634'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  -2.35%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.06270'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -2.24%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '20.77'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -4.74%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.07724'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -1.09%  '
$ws.Range('D12').Value = '1.660.18'
$ws.Range('E12').Value = '  -1.06%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '4.426'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  -1.85%  '
$ws.Range('D14').Value = '1.884.06'
$ws.Range('E14').Value = '  -1.23%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.5415'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -3.25%  '
$ws.Range('D16').Value = '0.0₅8108'
$ws.Range('E16').Value = '  -2.93%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '64.78'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -1.52%  '
$ws.Range('D18').Value = '26.156.13'
$ws.Range('E18').Value = '  -1.42%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '4.614'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -3.72%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '191.69'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -0.76%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '10.08'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -2.45%  '
$ws.Range('E23').Value = '  -4.84%  '
$ws.Range('E24').Value = '  +0.30%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '139.81'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +0.65%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '0.1224'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -3.72%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '7.182'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -3.12%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '16.07'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -1.44%  '
$ws.Range('E29').Value = '  -2.37%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.05974'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -5.46%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.271'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -1.68%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '3.548'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -1.62%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '3.254'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -5.49%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.599'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -5.45%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.9642'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -4.80%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '2.424'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -0.09%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '2.770'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -0.50%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.5691'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -7.78%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.01591'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -2.59%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '5.960'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -2.64%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.8548'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -0.99%  '
$ws.Range('E42').Value = '  +0.27%  '
$ws.Range('D43').Value = '1.005.92'
$ws.Range('E43').Value = '  -8.14%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '100.39'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -0.30%  '
$ws.Range('D45').Value = '1.798.66'
$ws.Range('E45').Value = '  -1.31%  '
$ws.Range('E46').Value = '  -1.12%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '56.67'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -3.38%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '1.006'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +1.17%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '7.995'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -2.51%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.05169'
$ws.Range('D50').Style = "Normal"
$ws.Range('B51').Value = 'RenderToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '1.447'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -3.18%  '
